$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.705.86'
$ws.Range("E2").Value = '  +2.84%  '
$ws.Range("D3").Value = '3.379.79'
$ws.Range("E3").Value = '  +3.67%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '191.23'
$ws.Range("E5").Value = '  +3.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '593.75'
$ws.Range("E6").Value = '  +2.20%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.609'
$ws.Range("E8").Value = '  +0.79%  '
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.75'
$ws.Range("E10").Value = '  +2.66%  '
$ws.Range("E11").Value = '  +1.74%  '
$ws.Range("D12").Value = '3.970.73'
$ws.Range("E12").Value = '  +3.82%  '
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.62'
$ws.Range("E14").Value = '  +3.18%  '
$ws.Range("D15").Value = '69.712.16'
$ws.Range("E15").Value = '  +2.91%  '
$ws.Range("E16").Value = '  +1.55%  '
$ws.Range("D17").Value = '3.395.22'
$ws.Range("E17").Value = '  +4.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '453.89'
$ws.Range("E18").Value = '  +15.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.82'
$ws.Range("E19").Value = '  +1.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.82'
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.79'
$ws.Range("E21").Value = '  +2.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.97'
$ws.Range("E22").Value = '  +6.19%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E25").Value = '  +3.20%  '
$ws.Range("E26").Value = '  +2.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.51'
$ws.Range("E27").Value = '  -0.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.02'
$ws.Range("E29").Value = '  +3.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.42'
$ws.Range("E30").Value = '  +3.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.60'
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("E32").Value = '  +2.54%  '
$ws.Range("E33").Value = '  +0.21%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  +6.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.50'
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.90'
$ws.Range("E38").Value = '  +4.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.814'
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("E41").Value = '  +2.09%  '
$ws.Range("D42").Value = '2.743.83'
$ws.Range("E42").Value = '  +4.98%  '
$ws.Range("E43").Value = '  +1.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.50'
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.14'
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '339.31'
$ws.Range("E47").Value = '  +1.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0285'
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.88'
$ws.Range("E50").Value = '  +4.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.33'
$ws.Range("E51").Value = '  -0.56%  '
